# Update the "取得日時" (retrieved datetime) column for the newly appended
# rows on the "ランサーズ" sheet: change the timestamp from
# 2026-02-01 12:43:33 to 2026-02-01 12:58:32 for rows 2 through 8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2026-02-01 12:43:33"
$newValue = "2026-02-01 12:58:32"

for ($row = 2; $row -le 8; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
